$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "59.375.93"
$ws.Range("E2").Value = "  +2.24%  "
Set-TextValue $ws.Range("D3") "2.593.70"
$ws.Range("E3").Value = "  +0.65%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.14%  "
Set-TextValue $ws.Range("D5") "531.45"
$ws.Range("E5").Value = "  +2.66%  "
Set-TextValue $ws.Range("D6") "140.54"
$ws.Range("E6").Value = "  +0.83%  "
Set-TextValue $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +1.05%  "
Set-TextValue $ws.Range("D9") "2.607.38"
$ws.Range("E9").Value = "  +0.69%  "
Set-TextValue $ws.Range("D10") "6.45"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  +2.18%  "
Set-TextValue $ws.Range("D12") "0.334"
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("E13").Value = "  +2.84%  "
Set-TextValue $ws.Range("D14") "3.051.30"
$ws.Range("E14").Value = "  +0.87%  "
Set-TextValue $ws.Range("D15") "59.116.45"
$ws.Range("E15").Value = "  +1.85%  "
Set-TextValue $ws.Range("D16") "20.45"
$ws.Range("E16").Value = "  +1.74%  "
Set-TextValue $ws.Range("D17") "0.0000134"
$ws.Range("E17").Value = "  +1.44%  "
Set-TextValue $ws.Range("D18") "2.598.59"
$ws.Range("E18").Value = "  +0.66%  "
Set-TextValue $ws.Range("D19") "346.95"
$ws.Range("E19").Value = "  +4.07%  "
Set-TextValue $ws.Range("D20") "4.33"
$ws.Range("E20").Value = "  +0.80%  "
Set-TextValue $ws.Range("D21") "10.11"
$ws.Range("E21").Value = "  +0.14%  "
Set-TextValue $ws.Range("D22") "6.39"
$ws.Range("E22").Value = "  +0.50%  "
Set-TextValue $ws.Range("D23") "0.999"
$ws.Range("E23").Value = "  -0.18%  "
Set-TextValue $ws.Range("D24") "67.50"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("E25").Value = "  +1.18%  "
Set-TextValue $ws.Range("D26") "0.406"
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("E27").Value = "  +0.33%  "
Set-TextValue $ws.Range("D28") "7.16"
$ws.Range("E28").Value = "  +3.11%  "
Set-TextValue $ws.Range("D29") "0.998"
$ws.Range("E29").Value = "  +0.06%  "
Set-TextValue $ws.Range("D30") "0.0₃0732"
$ws.Range("E30").Value = "  +1.86%  "
Set-TextValue $ws.Range("D31") "1.62"
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("E32").Value = "  -2.31%  "
Set-TextValue $ws.Range("D33") "18.78"
$ws.Range("E33").Value = "  +0.64%  "
Set-TextValue $ws.Range("D34") "149.81"
$ws.Range("E34").Value = "  +0.43%  "
Set-TextValue $ws.Range("D35") "3.97"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("E36").Value = "  +0.23%  "
Set-TextValue $ws.Range("D37") "36.86"
$ws.Range("E37").Value = "  +1.77%  "
Set-TextValue $ws.Range("D38") "1.48"
$ws.Range("E38").Value = "  +3.04%  "
Set-TextValue $ws.Range("D39") "0.835"
$ws.Range("E39").Value = "  -0.81%  "
Set-TextValue $ws.Range("D40") "0.828"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("E41").Value = "  +0.62%  "
Set-TextValue $ws.Range("D42") "0.999"
$ws.Range("E42").Value = "  +0.22%  "
Set-TextValue $ws.Range("D43") "271.98"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D44") "0.596"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D45") "10.74"
$ws.Range("E45").Value = "  +0.49%  "
Set-TextValue $ws.Range("D46") "0.0959"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D48") "4.64"
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D49") "1.947.08"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D50") "0.0221"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D51") "18.26"
$ws.Range("E51").Value = "  +1.64%  "
